# "Generate Report for Handoff"
# Adds two new handed-off files to the localization status report:
#   8567be4c-b30c-46c3-85f4-62ca48e51d69  (commit e27cbadae0e4305f524b21969f3a05bb2e472570)
#   ce8f7aa5-7343-4678-970f-1f723ba93e36  (commit a4a358ac64abc54f25868d9795d3946ad6dbdbbb)
# Appends one row per file to each of the three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$files = @(
    @{
        Guid       = "8567be4c-b30c-46c3-85f4-62ca48e51d69"
        CommitZh   = "e27cbadae0e4305f524b21969f3a05bb2e472570"
        CommitDe   = "e27cbadae0e4305f524b21969f3a05bb2e472570"
        OverviewDt = "2016-32-21 06:32:42"
        ZhHandoffDt= "2016-03-21 06:32:38"
        DeHandoffDt= "2016-03-21 06:32:42"
    },
    @{
        Guid       = "ce8f7aa5-7343-4678-970f-1f723ba93e36"
        CommitZh   = "a4a358ac64abc54f25868d9795d3946ad6dbdbbb"
        CommitDe   = "a4a358ac64abc54f25868d9795d3946ad6dbdbbb"
        OverviewDt = "2016-32-21 06:32:42"
        ZhHandoffDt= "2016-03-21 06:32:38"
        DeHandoffDt= "2016-03-21 06:32:42"
    }
)

$status = "Ready for handoff"
$ext = ".md"
$include = "Include"
$noDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Sheet: Overview  (columns: A=File Name, B=zh-cn, C=de-de, D=Latest Handoff Date)
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$startRow = 6

for ($i = 0; $i -lt $files.Count; $i++) {
    $f = $files[$i]
    $r = $startRow + $i
    $mdName = $f.Guid + $ext

    $wsOv.Range("A" + $r).Value = $mdName
    $wsOv.Range("B" + $r).Value = $status
    $wsOv.Range("C" + $r).Value = $status
    $wsOv.Range("D" + $r).Value = $f.OverviewDt

    $wsOv.Hyperlinks.Add($wsOv.Cells.Item($r, 1), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/" + $mdName, "", "", $mdName)
}

# ---------------------------------------------------------------------------
# Sheet: zh-cn  (A=Source File Name, B=File Extension, C=Status,
#                D=Latest Handoff File, E=Latest Handoff Datetime,
#                H=Handoff Reason, I=Dependency From)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$startRow = 6

for ($i = 0; $i -lt $files.Count; $i++) {
    $f = $files[$i]
    $r = $startRow + $i
    $mdName = $f.Guid + $ext
    $xlfName = $f.Guid + "." + $f.CommitZh + ".zh-cn.xlf"

    $wsZh.Range("A" + $r).Value = $mdName
    $wsZh.Range("B" + $r).Value = $ext
    $wsZh.Range("C" + $r).Value = $status
    $wsZh.Range("D" + $r).Value = $xlfName
    $wsZh.Range("E" + $r).Value = $f.ZhHandoffDt
    $wsZh.Range("E" + $r).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $wsZh.Range("H" + $r).Value = $noDate
    $wsZh.Range("I" + $r).Value = $include

    $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/" + $mdName
    $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $xlfName

    $wsZh.Hyperlinks.Add($wsZh.Cells.Item($r, 1), $mdUrl, "", "", $mdName)
    $wsZh.Hyperlinks.Add($wsZh.Cells.Item($r, 2), $mdUrl, "", "", $ext)
    $wsZh.Hyperlinks.Add($wsZh.Cells.Item($r, 4), $xlfUrl, "", "", $xlfName)
}

# ---------------------------------------------------------------------------
# Sheet: de-de  (same layout as zh-cn)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$startRow = 6

for ($i = 0; $i -lt $files.Count; $i++) {
    $f = $files[$i]
    $r = $startRow + $i
    $mdName = $f.Guid + $ext
    $xlfName = $f.Guid + "." + $f.CommitDe + ".de-de.xlf"

    $wsDe.Range("A" + $r).Value = $mdName
    $wsDe.Range("B" + $r).Value = $ext
    $wsDe.Range("C" + $r).Value = $status
    $wsDe.Range("D" + $r).Value = $xlfName
    $wsDe.Range("E" + $r).Value = $f.DeHandoffDt
    $wsDe.Range("E" + $r).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $wsDe.Range("H" + $r).Value = $noDate
    $wsDe.Range("I" + $r).Value = $include

    $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/" + $mdName
    $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $xlfName

    $wsDe.Hyperlinks.Add($wsDe.Cells.Item($r, 1), $mdUrl, "", "", $mdName)
    $wsDe.Hyperlinks.Add($wsDe.Cells.Item($r, 2), $mdUrl, "", "", $ext)
    $wsDe.Hyperlinks.Add($wsDe.Cells.Item($r, 4), $xlfUrl, "", "", $xlfName)
}

Write-Host "Handoff report rows added."
